# Update the "site_data" sheet with the new survey values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("site_data")

# river_length (row 7): 1400 -> 5000
$ws.Range("C7").Value = 5000

# area_catch (row 13): 12.4 -> 5.2
$ws.Range("C13").Value = 5.2

# Hq1pnat_catch (row 15): 0.107 -> 1.26
$ws.Range("C15").Value = 1.26

# area_urban_upstream (row 16): 4 -> 0
$ws.Range("C16").Value = 0

# Leave the cursor on the last-edited cell, matching the saved selection.
$ws.Activate()
$ws.Range("C15").Select()
